$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Merge the split "Mor" / "_GoBack" bookmark / "e robots..." runs
# back into a single run of plain text (removes the stray _GoBack bookmark).
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Mor" + "e robots can be purchased using the money raised.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "More robots can be purchased using the money raised.", 2)

# ---------------------------------------------------------------------------
# Change 2: Strike through the "Write Wagon Object creation code" bullet,
# matching the preceding bullet's formatting.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Write Wagon Object creation code*") {
        $p.Range.Font.StrikeThrough = $true
        break
    }
}

# ---------------------------------------------------------------------------
# Change 3: Add a new "WST_DESTROYED" row to the wagon-state table.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Cell(1, 1).Range.Text -like "State*") {
        $stateTable = $candidate
        break
    }
}
$newRow = $stateTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "WST_DESTROYED"
$notesCell = $newRow.Cells.Item(2).Range
$notesCell.Font.Italic = $false
$notesCell.Text = "No longer operating."

# ---------------------------------------------------------------------------
# Change 4: "Miner" -> "Driller" in the actionMaxCount description, leaving
# a _GoBack bookmark behind (mirrors Word's "last edit" bookmark).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Cell(1, 1).Range.Text -like "Instance*") {
        $instanceTable = $candidate
        break
    }
}
for ($r = 1; $r -le $instanceTable.Rows.Count; $r++) {
    if ($instanceTable.Cell($r, 1).Range.Text -like "actionMaxCount*") {
        $actionRow = $r
        break
    }
}
$cellRange = $instanceTable.Cell($actionRow, 3).Range
$rng = $d.Range($cellRange.Start, $cellRange.End)
$found = $rng.Find.Execute("Miner")
$rng.Text = "Driller"
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)
